# Auto-generated: applies the scheduled-runner price/profit refresh
# to each class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 221.76923
$ws.Range("I5").Value = 47.625
$ws.Range("J5").Value = 500.4
$ws.Range("K5").Value = 47.625
$ws.Range("L5").Value = 500.4
$ws.Range("M5").Value = 67.375
$ws.Range("N5").Value = -730.4

$ws.Range("H17").Value = 1409.2222
$ws.Range("J17").Value = 1481.3096
$ws.Range("L17").Value = 4443.9288
$ws.Range("N17").Value = -4779.9288

$ws.Range("H38").Value = 2017.2307
$ws.Range("I38").Value = 1845
$ws.Range("J38").Value = 2218.1667
$ws.Range("K38").Value = 5535
$ws.Range("L38").Value = 6654.500100000001
$ws.Range("M38").Value = -5163
$ws.Range("N38").Value = -7398.500100000001

$ws.Range("H106").Value = 4275612.5
$ws.Range("I106").Value = 4631705
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 4631705
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -4631074
$ws.Range("N106").Value = -3762

$ws.Range("H113").Value = 351001.66
$ws.Range("I113").Value = 351001.66
$ws.Range("K113").Value = 351001.66
$ws.Range("M113").Value = -347747.66

$ws.Range("H116").Value = 5324384.5
$ws.Range("J116").Value = 3249
$ws.Range("L116").Value = 3249
$ws.Range("N116").Value = -10133

$ws.Range("H137").Value = 32259122
$ws.Range("I137").Value = 34483764
$ws.Range("J137").Value = 1825
$ws.Range("K137").Value = 103451292
$ws.Range("L137").Value = 5475
$ws.Range("M137").Value = -103448742
$ws.Range("N137").Value = -10575

$ws.Range("H138").Value = 4632924.5
$ws.Range("I138").Value = 3477383.8
$ws.Range("J138").Value = 4833888
$ws.Range("K138").Value = 10432151.4
$ws.Range("L138").Value = 14501664
$ws.Range("M138").Value = -10427011.4
$ws.Range("N138").Value = -14511944

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 58833.332
$ws.Range("J131").Value = 58833.332
$ws.Range("L131").Value = 58833.332
$ws.Range("N131").Value = -68913.33199999999

$ws.Range("H132").Value = 2220.5
$ws.Range("I132").Value = 1580.6897
$ws.Range("J132").Value = 4282.1113
$ws.Range("K132").Value = 4742.0691
$ws.Range("L132").Value = 12846.3339
$ws.Range("M132").Value = -2212.0691
$ws.Range("N132").Value = -17906.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 30666.666
$ws.Range("I5").Value = 26000
$ws.Range("J5").Value = 40000
$ws.Range("K5").Value = 26000
$ws.Range("L5").Value = 40000
$ws.Range("M5").Value = -25887
$ws.Range("N5").Value = -40226

$ws.Range("H8").Value = 4052.5
$ws.Range("I8").Value = 105
$ws.Range("K8").Value = 105
$ws.Range("M8").Value = 35

$ws.Range("H105").Value = 3369.577
$ws.Range("I105").Value = 3200.5
$ws.Range("J105").Value = 3750
$ws.Range("K105").Value = 3200.5
$ws.Range("L105").Value = 3750
$ws.Range("M105").Value = -1453.5
$ws.Range("N105").Value = -7244

$ws.Range("H107").Value = 916.2
$ws.Range("I107").Value = 916.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 916.2
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = 1003.8
$ws.Range("M107").ClearContents()

$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1301.9259
$ws.Range("I31").Value = 855.381
$ws.Range("J31").Value = 2864.8333
$ws.Range("K31").Value = 855.381
$ws.Range("L31").Value = 2864.8333
$ws.Range("M31").Value = -560.381
$ws.Range("N31").Value = -3454.8333

$ws.Range("H34").Value = 1301.9259
$ws.Range("I34").Value = 855.381
$ws.Range("J34").Value = 2864.8333
$ws.Range("K34").Value = 855.381
$ws.Range("L34").Value = 2864.8333
$ws.Range("M34").Value = -653.381
$ws.Range("N34").Value = -3268.8333

$ws.Range("H58").Value = 4166.7144
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4166.7144
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = 4166.7144
$ws.Range("N58").Value = -4572.7144
$ws.Range("L58").ClearContents()

$ws.Range("H99").Value = 6945851
$ws.Range("I99").Value = 10417977
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 10417977
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -10416479
$ws.Range("N99").Value = -4596

$ws.Range("H105").Value = 788.2778
$ws.Range("I105").Value = 786.8125
$ws.Range("K105").Value = 786.8125
$ws.Range("M105").Value = 960.1875

$ws.Range("H117").Value = 27356
$ws.Range("J117").Value = 27356
$ws.Range("L117").Value = 27356
$ws.Range("N117").Value = -36534

$ws.Range("H122").Value = 1581.1538
$ws.Range("I122").Value = 879.4737
$ws.Range("J122").Value = 3485.7144
$ws.Range("K122").Value = 2638.4211
$ws.Range("L122").Value = 10457.1432
$ws.Range("M122").Value = -188.4211
$ws.Range("N122").Value = -15357.1432

$ws.Range("H126").Value = 6945851
$ws.Range("I126").Value = 10417977
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 31253931
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -31251461
$ws.Range("N126").Value = -9740

$ws.Range("H136").Value = 4166.7144
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4166.7144
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = 12500.1432
$ws.Range("N136").Value = -17600.1432
$ws.Range("L136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3505.5
$ws.Range("I70").Value = 2011
$ws.Range("K70").Value = 6033
$ws.Range("M70").Value = -5718

$ws.Range("H73").Value = 3505.5
$ws.Range("I73").Value = 2011
$ws.Range("K73").Value = 6033
$ws.Range("M73").Value = -4941

$ws.Range("H80").Value = 1120

$ws.Range("H83").Value = 1120

$ws.Range("H140").Value = 4054.366
$ws.Range("I140").Value = 5036.696
$ws.Range("J140").Value = 2799.1667
$ws.Range("K140").Value = 15110.088
$ws.Range("L140").Value = 8397.500100000001
$ws.Range("M140").Value = -9930.088
$ws.Range("N140").Value = -18757.5001

$ws.Range("H141").Value = 3923.7693
$ws.Range("I141").Value = 4876.125
$ws.Range("J141").Value = 2400
$ws.Range("K141").Value = 14628.375
$ws.Range("L141").Value = 7200
$ws.Range("M141").Value = -9448.375
$ws.Range("N141").Value = -17560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1112269.5
$ws.Range("I122").Value = 1235757.2
$ws.Range("J122").Value = 880
$ws.Range("K122").Value = 3707271.6
$ws.Range("L122").Value = 2640
$ws.Range("M122").Value = -3704821.6
$ws.Range("N122").Value = -7540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1783.5555
$ws.Range("I46").Value = 1100.1111
$ws.Range("J46").Value = 2467
$ws.Range("K46").Value = 1100.1111
$ws.Range("L46").Value = 2467
$ws.Range("M46").Value = -912.1111000000001
$ws.Range("N46").Value = -2843

$ws.Range("H132").Value = 4836.0454
$ws.Range("I132").Value = 4226.5
$ws.Range("J132").Value = 5567.5
$ws.Range("K132").Value = 12679.5
$ws.Range("L132").Value = 16702.5
$ws.Range("M132").Value = -10149.5
$ws.Range("N132").Value = -21762.5

$ws.Range("H136").Value = 3109.0544
$ws.Range("I136").Value = 1776.4681
$ws.Range("K136").Value = 5329.4043
$ws.Range("M136").Value = -2779.4043

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 29000
$ws.Range("J121").Value = 29000
$ws.Range("L121").Value = 29000
$ws.Range("N121").Value = -32494

$ws.Range("H123").Value = 27913.637
$ws.Range("J123").Value = 27913.637
$ws.Range("L123").Value = 27913.637
$ws.Range("N123").Value = -37713.637

$ws.Range("H132").Value = 13160072
$ws.Range("I132").Value = 16668337
$ws.Range("J132").Value = 4081
$ws.Range("K132").Value = 50005011
$ws.Range("L132").Value = 12243
$ws.Range("M132").Value = -50002481
$ws.Range("N132").Value = -17303
